{"js": "const replacements = [\n  [\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"],\n  [\"19\u00f78=\", \"49\u00f75=\"],\n  [\"45\u00f78=\", \"13\u00f78=\"],\n  [\"85\u00f76=\", \"90\u00f77=\"],\n  [\"36\u00f75=\", \"95\u00f78=\"],\n  [\"66\u00f79=\", \"87\u00f73=\"],\n  [\"87\u00f78=\", \"81\u00f73=\"],\n  [\"22\u00f72=\", \"64\u00f75=\"],\n  [\"27\u00f78=\", \"68\u00f79=\"],\n  [\"28\u00f72=\", \"38\u00f76=\"],\n  [\"21\u00f78=\", \"23\u00f75=\"],\n  [\"59\u00f72=\", \"98\u00f72=\"],\n  [\"34\u00f75=\", \"93\u00f74=\"],\n  [\"23\u00f77=\", \"36\u00f72=\"],\n  [\"22\u00f77=\", \"89\u00f77=\"],\n  [\"43\u00f79=\", \"27\u00f79=\"],\n  [\"99\u00f75=\", \"83\u00f78=\"],\n  [\"55\u00f76=\", \"44\u00f72=\"],\n  [\"30\u00f79=\", \"29\u00f76=\"],\n  [\"25\u00f79=\", \"65\u00f73=\"],\n  [\"38\u00f79=\", \"10\u00f72=\"],\n  [\"35\u00f73=\", \"26\u00f75=\"],\n  [\"67\u00f79=\", \"10\u00f73=\"],\n  [\"64\u00f73=\", \"22\u00f73=\"],\n  [\"66\u00f74=\", \"54\u00f79=\"],\n  [\"98\u00f75=\", \"72\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"),\n    @(\"19\u00f78=\", \"49\u00f75=\"),\n    @(\"45\u00f78=\", \"13\u00f78=\"),\n    @(\"85\u00f76=\", \"90\u00f77=\"),\n    @(\"36\u00f75=\", \"95\u00f78=\"),\n    @(\"66\u00f79=\", \"87\u00f73=\"),\n    @(\"87\u00f78=\", \"81\u00f73=\"),\n    @(\"22\u00f72=\", \"64\u00f75=\"),\n    @(\"27\u00f78=\", \"68\u00f79=\"),\n    @(\"28\u00f72=\", \"38\u00f76=\"),\n    @(\"21\u00f78=\", \"23\u00f75=\"),\n    @(\"59\u00f72=\", \"98\u00f72=\"),\n    @(\"34\u00f75=\", \"93\u00f74=\"),\n    @(\"23\u00f77=\", \"36\u00f72=\"),\n    @(\"22\u00f77=\", \"89\u00f77=\"),\n    @(\"43\u00f79=\", \"27\u00f79=\"),\n    @(\"99\u00f75=\", \"83\u00f78=\"),\n    @(\"55\u00f76=\", \"44\u00f72=\"),\n    @(\"30\u00f79=\", \"29\u00f76=\"),\n    @(\"25\u00f79=\", \"65\u00f73=\"),\n    @(\"38\u00f79=\", \"10\u00f72=\"),\n    @(\"35\u00f73=\", \"26\u00f75=\"),\n    @(\"67\u00f79=\", \"10\u00f73=\"),\n    @(\"64\u00f73=\", \"22\u00f73=\"),\n    @(\"66\u00f74=\", \"54\u00f79=\"),\n    @(\"98\u00f75=\", \"72\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
